# Refatoração no caminho que está genérico, criação de um parametro para
# saber se o teste passou ou não e ajustes no screenshot

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cadastro")

# New "Resultado" column with a pass/fail flag
$ws.Range("M1").Value = "Resultado"
$ws.Range("M3").Value = "Reprovado"

# Username column: FLP00 -> BRUN227
$ws.Range("A2").Value = "BRUN227"
$ws.Range("A3").Value = "BRUN227"

$ws.Range("M2").Value = "Aprovado"

# Match header style (bold + centered) used by the rest of row 1
$ws.Range("M1").Font.Bold = $true
$ws.Range("M1").HorizontalAlignment = -4108

# M2/M3 keep the plain default style (not the centered column default)
$ws.Range("M2").Style = "Normal"
$ws.Range("M3").Style = "Normal"

$ws.Columns.Item(13).ColumnWidth = 9.7

$ws.Range("A3").Select()

$wb.Save()
